$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 3440.9167  # H32: 3790.6956 -> 3440.9167
$ws.Cells.Item(32, 9).Value = 3543  # I32: 3693.5 -> 3543
$ws.Cells.Item(32, 10).Value = 3270.7778  # J32: 3941.889 -> 3270.7778
$ws.Cells.Item(32, 11).Value = 3543  # K32: 3693.5 -> 3543
$ws.Cells.Item(32, 12).Value = 3270.7778  # L32: 3941.889 -> 3270.7778
$ws.Cells.Item(32, 13).Value = -3217  # M32: -3367.5 -> -3217
$ws.Cells.Item(32, 14).Value = -3922.7778  # N32: -4593.889 -> -3922.7778
$ws.Cells.Item(87, 8).Value = 113332.336  # H87: 123325.664 -> 113332.336
$ws.Cells.Item(87, 10).Value = 113332.336  # J87: 123325.664 -> 113332.336
$ws.Cells.Item(87, 12).Value = 113332.336  # L87: 123325.664 -> 113332.336
$ws.Cells.Item(87, 14).Value = -115828.336  # N87: -125821.664 -> -115828.336
$ws.Cells.Item(90, 8).Value = 113332.336  # H90: 123325.664 -> 113332.336
$ws.Cells.Item(90, 10).Value = 113332.336  # J90: 123325.664 -> 113332.336
$ws.Cells.Item(90, 12).Value = 339997.008  # L90: 369976.992 -> 339997.008
$ws.Cells.Item(90, 14).Value = -352477.008  # N90: -382456.992 -> -352477.008
$ws.Cells.Item(94, 8).Value = 50130960  # H94: 55700812 -> 50130960
$ws.Cells.Item(94, 9).Value = 83344100  # I94: 100012456 -> 83344100
$ws.Cells.Item(94, 11).Value = 83344100  # K94: 100012456 -> 83344100
$ws.Cells.Item(94, 13).Value = -83343649  # M94: -100012005 -> -83343649
$ws.Cells.Item(112, 8).Value = 49210.715  # H112: 64036.562 -> 49210.715
$ws.Cells.Item(112, 10).Value = 85173.75  # J112: 144749.28 -> 85173.75
$ws.Cells.Item(112, 12).Value = 255521.25  # L112: 434247.84 -> 255521.25
$ws.Cells.Item(112, 14).Value = -257737.25  # N112: -436463.84 -> -257737.25
$ws.Cells.Item(113, 8).Value = 11948.833  # H113: 14770.5 -> 11948.833
$ws.Cells.Item(113, 9).Value = 14759.75  # I113: 17578 -> 14759.75
$ws.Cells.Item(113, 10).Value = 6327  # J113: 6348 -> 6327
$ws.Cells.Item(113, 11).Value = 14759.75  # K113: 17578 -> 14759.75
$ws.Cells.Item(113, 12).Value = 6327  # L113: 6348 -> 6327
$ws.Cells.Item(113, 13).Value = -11505.75  # M113: -14324 -> -11505.75
$ws.Cells.Item(113, 14).Value = -12835  # N113: -12856 -> -12835
$ws.Cells.Item(116, 8).Value = 12226286  # H116: 1393781 -> 12226286
$ws.Cells.Item(116, 9).Value = 17462268  # I116: 2224652 -> 17462268
$ws.Cells.Item(116, 11).Value = 17462268  # K116: 2224652 -> 17462268
$ws.Cells.Item(116, 13).Value = -17458826  # M116: -2221210 -> -17458826

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5637.923  # H2: 5456.8887 -> 5637.923
$ws.Cells.Item(2, 9).Value = 6490.5  # I2: 6217.143 -> 6490.5
$ws.Cells.Item(2, 11).Value = 6490.5  # K2: 6217.143 -> 6490.5
$ws.Cells.Item(2, 13).Value = -6377.5  # M2: -6104.143 -> -6377.5
$ws.Cells.Item(45, 8).Value = 8544  # H45: 7353.636 -> 8544
$ws.Cells.Item(45, 9).Value = 7779.4  # I45: 6127.2856 -> 7779.4
$ws.Cells.Item(45, 11).Value = 7779.4  # K45: 6127.2856 -> 7779.4
$ws.Cells.Item(45, 13).Value = -7402.4  # M45: -5750.2856 -> -7402.4
$ws.Cells.Item(61, 8).Value = 3889.5264  # H61: 3973.5676 -> 3889.5264
$ws.Cells.Item(61, 9).Value = 3262.077  # I61: 3361.36 -> 3262.077
$ws.Cells.Item(61, 11).Value = 3262.077  # K61: 3361.36 -> 3262.077
$ws.Cells.Item(61, 13).Value = -3050.077  # M61: -3149.36 -> -3050.077
$ws.Cells.Item(88, 8).Value = 2747  # H88: 2785 -> 2747
$ws.Cells.Item(88, 10).Value = 4256.857  # J88: 4959.6 -> 4256.857
$ws.Cells.Item(88, 12).Value = 4256.857  # L88: 4959.6 -> 4256.857
$ws.Cells.Item(88, 14).Value = -5068.857  # N88: -5771.6 -> -5068.857
$ws.Cells.Item(91, 8).Value = 2747  # H91: 2785 -> 2747
$ws.Cells.Item(91, 10).Value = 4256.857  # J91: 4959.6 -> 4256.857
$ws.Cells.Item(91, 12).Value = 4256.857  # L91: 4959.6 -> 4256.857
$ws.Cells.Item(91, 14).Value = -7064.857  # N91: -7767.6 -> -7064.857
$ws.Cells.Item(101, 8).Value = 79996  # H101: 79996.5 -> 79996
$ws.Cells.Item(101, 10).Value = 79996  # J101: 79996.5 -> 79996
$ws.Cells.Item(101, 12).Value = 79996  # L101: 79996.5 -> 79996
$ws.Cells.Item(101, 14).Value = -86486  # N101: -86486.5 -> -86486
$ws.Cells.Item(116, 8).Value = 5637.923  # H116: 5456.8887 -> 5637.923
$ws.Cells.Item(116, 9).Value = 6490.5  # I116: 6217.143 -> 6490.5
$ws.Cells.Item(116, 11).Value = 6490.5  # K116: 6217.143 -> 6490.5
$ws.Cells.Item(116, 13).Value = -4196.5  # M116: -3923.143 -> -4196.5
$ws.Cells.Item(124, 8).Value = 39429  # H124: 0 -> 39429
$ws.Cells.Item(124, 10).Value = 39429  # J124: 0 -> 39429
$ws.Cells.Item(124, 12).Value = 39429  # L124: 0 -> 39429
$ws.Cells.Item(124, 14).Value = -49249  # N124: None -> -49249
$ws.Cells.Item(136, 8).Value = 3889.5264  # H136: 3973.5676 -> 3889.5264
$ws.Cells.Item(136, 9).Value = 3262.077  # I136: 3361.36 -> 3262.077
$ws.Cells.Item(136, 11).Value = 9786.231  # K136: 10084.08 -> 9786.231
$ws.Cells.Item(136, 13).Value = -7236.231  # M136: -7534.08 -> -7236.231

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5637.923  # H3: 5456.8887 -> 5637.923
$ws.Cells.Item(3, 9).Value = 6490.5  # I3: 6217.143 -> 6490.5
$ws.Cells.Item(3, 11).Value = 6490.5  # K3: 6217.143 -> 6490.5
$ws.Cells.Item(3, 13).Value = -6376.5  # M3: -6103.143 -> -6376.5
$ws.Cells.Item(22, 8).Value = 598.4286  # H22: 614.8333 -> 598.4286
$ws.Cells.Item(22, 9).Value = 598.4286  # I22: 614.8333 -> 598.4286
$ws.Cells.Item(22, 11).Value = 598.4286  # K22: 614.8333 -> 598.4286
$ws.Cells.Item(22, 13).Value = -425.4286  # M22: -441.8333 -> -425.4286
$ws.Cells.Item(88, 8).Value = 25838.285  # H88: 28478.166 -> 25838.285
$ws.Cells.Item(88, 9).Value = 10000  # I88: 0 -> 10000
$ws.Cells.Item(88, 10).Value = 28478  # J88: 28478.166 -> 28478
$ws.Cells.Item(88, 11).Value = 10000  # K88: 0 -> 10000
$ws.Cells.Item(88, 12).Value = 28478  # L88: 28478.166 -> 28478
$ws.Cells.Item(88, 13).Value = -9594  # M88: None -> -9594
$ws.Cells.Item(88, 14).Value = -29290  # N88: -29290.166 -> -29290
$ws.Cells.Item(91, 8).Value = 25838.285  # H91: 28478.166 -> 25838.285
$ws.Cells.Item(91, 9).Value = 10000  # I91: 0 -> 10000
$ws.Cells.Item(91, 10).Value = 28478  # J91: 28478.166 -> 28478
$ws.Cells.Item(91, 11).Value = 10000  # K91: 0 -> 10000
$ws.Cells.Item(91, 12).Value = 28478  # L91: 28478.166 -> 28478
$ws.Cells.Item(91, 13).Value = -8596  # M91: None -> -8596
$ws.Cells.Item(91, 14).Value = -31286  # N91: -31286.166 -> -31286
$ws.Cells.Item(109, 8).Value = 0  # H109: 59999 -> 0
$ws.Cells.Item(109, 10).Value = 0  # J109: 59999 -> 0
$ws.Cells.Item(109, 12).Value = 0  # L109: 59999 -> 0
$ws.Cells.Item(109, 14).ClearContents()  # N109: -62773 -> (removed)
$ws.Cells.Item(130, 8).Value = 77854.5  # H130: 120000 -> 77854.5
$ws.Cells.Item(130, 9).Value = 35709  # I130: 0 -> 35709
$ws.Cells.Item(130, 11).Value = 35709  # K130: 0 -> 35709
$ws.Cells.Item(130, 13).Value = -30689  # M130: None -> -30689

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1889.8  # H16: 1895 -> 1889.8
$ws.Cells.Item(16, 9).Value = 1658.3334  # I16: 1618.75 -> 1658.3334
$ws.Cells.Item(16, 10).Value = 2237  # J16: 3000 -> 2237
$ws.Cells.Item(16, 11).Value = 1658.3334  # K16: 1618.75 -> 1658.3334
$ws.Cells.Item(16, 12).Value = 2237  # L16: 3000 -> 2237
$ws.Cells.Item(16, 13).Value = -1371.3334  # M16: -1331.75 -> -1371.3334
$ws.Cells.Item(16, 14).Value = -2811  # N16: -3574 -> -2811
$ws.Cells.Item(43, 8).Value = 11887.4  # H43: 12067.5 -> 11887.4
$ws.Cells.Item(43, 10).Value = 11887.4  # J43: 12067.5 -> 11887.4
$ws.Cells.Item(43, 12).Value = 11887.4  # L43: 12067.5 -> 11887.4
$ws.Cells.Item(43, 14).Value = -12255.4  # N43: -12435.5 -> -12255.4
$ws.Cells.Item(68, 8).Value = 39999.668  # H68: 59666.668 -> 39999.668
$ws.Cells.Item(68, 10).Value = 39999.668  # J68: 59666.668 -> 39999.668
$ws.Cells.Item(68, 12).Value = 39999.668  # L68: 59666.668 -> 39999.668
$ws.Cells.Item(68, 14).Value = -41497.668  # N68: -61164.668 -> -41497.668
$ws.Cells.Item(71, 8).Value = 39999.668  # H71: 59666.668 -> 39999.668
$ws.Cells.Item(71, 10).Value = 39999.668  # J71: 59666.668 -> 39999.668
$ws.Cells.Item(71, 12).Value = 119999.004  # L71: 179000.004 -> 119999.004
$ws.Cells.Item(71, 14).Value = -127487.004  # N71: -186488.004 -> -127487.004
$ws.Cells.Item(100, 8).Value = 95996.664  # H100: 99330.336 -> 95996.664
$ws.Cells.Item(100, 10).Value = 95996.664  # J100: 99330.336 -> 95996.664
$ws.Cells.Item(100, 12).Value = 95996.664  # L100: 99330.336 -> 95996.664
$ws.Cells.Item(100, 14).Value = -98160.664  # N100: -101494.336 -> -98160.664
$ws.Cells.Item(101, 8).Value = 11887.4  # H101: 12067.5 -> 11887.4
$ws.Cells.Item(101, 10).Value = 11887.4  # J101: 12067.5 -> 11887.4
$ws.Cells.Item(101, 12).Value = 11887.4  # L101: 12067.5 -> 11887.4
$ws.Cells.Item(101, 14).Value = -18377.4  # N101: -18557.5 -> -18377.4
$ws.Cells.Item(113, 8).Value = 1889.8  # H113: 1895 -> 1889.8
$ws.Cells.Item(113, 9).Value = 1658.3334  # I113: 1618.75 -> 1658.3334
$ws.Cells.Item(113, 10).Value = 2237  # J113: 3000 -> 2237
$ws.Cells.Item(113, 11).Value = 1658.3334  # K113: 1618.75 -> 1658.3334
$ws.Cells.Item(113, 12).Value = 2237  # L113: 3000 -> 2237
$ws.Cells.Item(113, 13).Value = 511.6666  # M113: 551.25 -> 511.6666
$ws.Cells.Item(113, 14).Value = -6577  # N113: -7340 -> -6577

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(41, 8).Value = 2155.4583  # H41: 2161.2917 -> 2155.4583
$ws.Cells.Item(41, 9).Value = 4047.889  # I41: 4546.375 -> 4047.889
$ws.Cells.Item(41, 10).Value = 1020  # J41: 968.75 -> 1020
$ws.Cells.Item(41, 11).Value = 12143.667  # K41: 13639.125 -> 12143.667
$ws.Cells.Item(41, 12).Value = 3060  # L41: 2906.25 -> 3060
$ws.Cells.Item(41, 13).Value = -11805.667  # M41: -13301.125 -> -11805.667
$ws.Cells.Item(41, 14).Value = -3736  # N41: -3582.25 -> -3736
$ws.Cells.Item(51, 8).Value = 1823.75  # H51: 1909.1052 -> 1823.75
$ws.Cells.Item(51, 9).Value = 756.3077  # I51: 802.5 -> 756.3077
$ws.Cells.Item(51, 11).Value = 2268.9231  # K51: 2407.5 -> 2268.9231
$ws.Cells.Item(51, 13).Value = -1808.9231  # M51: -1947.5 -> -1808.9231
$ws.Cells.Item(68, 8).Value = 11716.333  # H68: 12330.909 -> 11716.333
$ws.Cells.Item(68, 9).Value = 405  # I68: 409.75 -> 405
$ws.Cells.Item(68, 10).Value = 17372  # J68: 19143 -> 17372
$ws.Cells.Item(68, 11).Value = 1215  # K68: 1229.25 -> 1215
$ws.Cells.Item(68, 12).Value = 52116  # L68: 57429 -> 52116
$ws.Cells.Item(68, 13).Value = -404  # M68: -418.25 -> -404
$ws.Cells.Item(68, 14).Value = -53738  # N68: -59051 -> -53738
$ws.Cells.Item(71, 8).Value = 11716.333  # H71: 12330.909 -> 11716.333
$ws.Cells.Item(71, 9).Value = 405  # I71: 409.75 -> 405
$ws.Cells.Item(71, 10).Value = 17372  # J71: 19143 -> 17372
$ws.Cells.Item(71, 11).Value = 3645  # K71: 3687.75 -> 3645
$ws.Cells.Item(71, 12).Value = 156348  # L71: 172287 -> 156348
$ws.Cells.Item(71, 13).Value = 411  # M71: 368.25 -> 411
$ws.Cells.Item(71, 14).Value = -164460  # N71: -180399 -> -164460
$ws.Cells.Item(107, 8).Value = 1244  # H107: 1309.25 -> 1244
$ws.Cells.Item(107, 10).Value = 1818  # J107: 1979.8 -> 1818
$ws.Cells.Item(107, 12).Value = 5454  # L107: 5939.4 -> 5454
$ws.Cells.Item(107, 14).Value = -9294  # N107: -9779.4 -> -9294
$ws.Cells.Item(140, 8).Value = 13041.533  # H140: 13997.714 -> 13041.533
$ws.Cells.Item(140, 9).Value = 13041.533  # I140: 14997.538 -> 13041.533
$ws.Cells.Item(140, 10).Value = 0  # J140: 1000 -> 0
$ws.Cells.Item(140, 11).Value = 39124.599  # K140: 44992.614 -> 39124.599
$ws.Cells.Item(140, 12).Value = 0  # L140: 3000 -> 0
$ws.Cells.Item(140, 13).Value = -33944.599  # M140: -39812.614 -> -33944.599
$ws.Cells.Item(140, 14).ClearContents()  # N140: -13360 -> (removed)

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62, 8).Value = 17000  # H62: 38000 -> 17000
$ws.Cells.Item(62, 9).Value = 17000  # I62: 48000 -> 17000
$ws.Cells.Item(62, 10).Value = 0  # J62: 28000 -> 0
$ws.Cells.Item(62, 11).Value = 17000  # K62: 48000 -> 17000
$ws.Cells.Item(62, 12).Value = 0  # L62: 28000 -> 0
$ws.Cells.Item(62, 13).Value = -16314  # M62: -47314 -> -16314
$ws.Cells.Item(62, 14).ClearContents()  # N62: -29372 -> (removed)
$ws.Cells.Item(65, 8).Value = 17000  # H65: 38000 -> 17000
$ws.Cells.Item(65, 9).Value = 17000  # I65: 48000 -> 17000
$ws.Cells.Item(65, 10).Value = 0  # J65: 28000 -> 0
$ws.Cells.Item(65, 11).Value = 51000  # K65: 144000 -> 51000
$ws.Cells.Item(65, 12).Value = 0  # L65: 84000 -> 0
$ws.Cells.Item(65, 13).Value = -47568  # M65: -140568 -> -47568
$ws.Cells.Item(65, 14).ClearContents()  # N65: -90864 -> (removed)
$ws.Cells.Item(70, 8).Value = 5989  # H70: 4494 -> 5989
$ws.Cells.Item(70, 9).Value = 5989  # I70: 8 -> 5989
$ws.Cells.Item(70, 10).Value = 5989  # J70: 5989.3335 -> 5989
$ws.Cells.Item(70, 11).Value = 5989  # K70: 8 -> 5989
$ws.Cells.Item(70, 12).Value = 5989  # L70: 5989.3335 -> 5989
$ws.Cells.Item(70, 13).Value = -5719  # M70: 262 -> -5719
$ws.Cells.Item(70, 14).Value = -6529  # N70: -6529.3335 -> -6529
$ws.Cells.Item(73, 8).Value = 5989  # H73: 4494 -> 5989
$ws.Cells.Item(73, 9).Value = 5989  # I73: 8 -> 5989
$ws.Cells.Item(73, 10).Value = 5989  # J73: 5989.3335 -> 5989
$ws.Cells.Item(73, 11).Value = 5989  # K73: 8 -> 5989
$ws.Cells.Item(73, 12).Value = 5989  # L73: 5989.3335 -> 5989
$ws.Cells.Item(73, 13).Value = -5053  # M73: 928 -> -5053
$ws.Cells.Item(73, 14).Value = -7861  # N73: -7861.3335 -> -7861
$ws.Cells.Item(92, 8).Value = 12928.857  # H92: 14050.2 -> 12928.857
$ws.Cells.Item(92, 10).Value = 12928.857  # J92: 14050.2 -> 12928.857
$ws.Cells.Item(92, 12).Value = 12928.857  # L92: 14050.2 -> 12928.857
$ws.Cells.Item(92, 14).Value = -16672.857  # N92: -17794.2 -> -16672.857

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(63, 8).Value = 35000  # H63: 55000 -> 35000
$ws.Cells.Item(63, 9).Value = 35000  # I63: 55000 -> 35000
$ws.Cells.Item(63, 11).Value = 35000  # K63: 55000 -> 35000
$ws.Cells.Item(63, 13).Value = -34251  # M63: -54251 -> -34251
$ws.Cells.Item(66, 8).Value = 35000  # H66: 55000 -> 35000
$ws.Cells.Item(66, 9).Value = 35000  # I66: 55000 -> 35000
$ws.Cells.Item(66, 11).Value = 105000  # K66: 165000 -> 105000
$ws.Cells.Item(66, 13).Value = -101256  # M66: -161256 -> -101256
$ws.Cells.Item(101, 8).Value = 32707.125  # H101: 32707.25 -> 32707.125
$ws.Cells.Item(101, 10).Value = 32707.125  # J101: 32707.25 -> 32707.125
$ws.Cells.Item(101, 12).Value = 32707.125  # L101: 32707.25 -> 32707.125
$ws.Cells.Item(101, 14).Value = -39197.125  # N101: -39197.25 -> -39197.125
$ws.Cells.Item(104, 8).Value = 21665  # H104: 22960 -> 21665
$ws.Cells.Item(104, 10).Value = 21665  # J104: 22960 -> 21665
$ws.Cells.Item(104, 12).Value = 21665  # L104: 22960 -> 21665
$ws.Cells.Item(104, 14).Value = -28653  # N104: -29948 -> -28653
$ws.Cells.Item(132, 8).Value = 624294.9399999999  # H132: 748261.9 -> 624294.9399999999
$ws.Cells.Item(132, 9).Value = 830815.4399999999  # I132: 996362.75 -> 830815.4399999999
$ws.Cells.Item(132, 10).Value = 4733.3335  # J132: 3959.4 -> 4733.3335
$ws.Cells.Item(132, 11).Value = 2492446.32  # K132: 2989088.25 -> 2492446.32
$ws.Cells.Item(132, 12).Value = 14200.0005  # L132: 11878.2 -> 14200.0005
$ws.Cells.Item(132, 13).Value = -2489916.32  # M132: -2986558.25 -> -2489916.32
$ws.Cells.Item(132, 14).Value = -19260.0005  # N132: -16938.2 -> -19260.0005
$ws.Cells.Item(136, 8).Value = 4677.3213  # H136: 4840.625 -> 4677.3213
$ws.Cells.Item(136, 9).Value = 1865.1333  # I136: 1992.0769 -> 1865.1333
$ws.Cells.Item(136, 10).Value = 7922.154  # J136: 8207.091 -> 7922.154
$ws.Cells.Item(136, 11).Value = 5595.3999  # K136: 5976.2307 -> 5595.3999
$ws.Cells.Item(136, 12).Value = 23766.462  # L136: 24621.273 -> 23766.462
$ws.Cells.Item(136, 13).Value = -3045.3999  # M136: -3426.2307 -> -3045.3999
$ws.Cells.Item(136, 14).Value = -28866.462  # N136: -29721.273 -> -28866.462

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 15475.714  # H2: 30000 -> 15475.714
$ws.Cells.Item(2, 9).Value = 15475.714  # I2: 30000 -> 15475.714
$ws.Cells.Item(2, 11).Value = 15475.714  # K2: 30000 -> 15475.714
$ws.Cells.Item(2, 13).Value = -15363.714  # M2: -29888 -> -15363.714
$ws.Cells.Item(105, 8).Value = 49663.668  # H105: 50000 -> 49663.668
$ws.Cells.Item(105, 10).Value = 49663.668  # J105: 50000 -> 49663.668
$ws.Cells.Item(105, 12).Value = 49663.668  # L105: 50000 -> 49663.668
$ws.Cells.Item(105, 14).Value = -56651.668  # N105: -56988 -> -56651.668
